# Fix minutes/seconds formatting (zero-pad to 2 digits) in the "Общее время"
# (haul) column of the sprint-top-by-haul stats sheet.
#
# Values look like "87 ч. 45 мин. 3 сек." and must become
# "87 ч. 45 мин. 03 сек." — i.e. the hours stay unpadded, while the minutes
# and seconds components are zero-padded to (at least) two digits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRows = $ws.UsedRange.Rows.Count

# Column D holds "Общее время" (the haul time column) per the header row.
$rng = $ws.Range("D2:D$usedRows")
$vals = $rng.Value2

$regex = [regex]'^(\d+)\s*ч\.\s*(\d+)\s*мин\.\s*(\d+)\s*сек\.$'

$rows = $vals.GetLength(0)
$changed = 0

for ($i = 1; $i -le $rows; $i++) {
    $cellValue = $vals[$i, 1]
    if ($cellValue -ne $null) {
        $text = [string]$cellValue
        $m = $regex.Match($text)
        if ($m.Success) {
            $hours = $m.Groups[1].Value
            $minutes = $m.Groups[2].Value
            $seconds = $m.Groups[3].Value

            $minutesPadded = $minutes.PadLeft(2, '0')
            $secondsPadded = $seconds.PadLeft(2, '0')

            $newText = "$hours ч. $minutesPadded мин. $secondsPadded сек."

            if (-not $newText.Equals($text)) {
                $vals[$i, 1] = $newText
                $changed = $changed + 1
            }
        }
    }
}

$rng.Value2 = $vals

Write-Host "Updated $changed haul time value(s) in column D."
